$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header labels: "<Name>_old" -> "<Name>_FV2404" and
#    "<Name>_new" -> "<Name>_FV2410" (the "diff" header stays as-is).
# ---------------------------------------------------------------------------
$fv2404Headers = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2404Headers[$i]
}

for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2410Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into a real Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U84"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split/freeze after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
